# Applies the OOXML diff described in the commit:
#   "fixed isLeaf() to count children and compare to number of threaded children"
#
# Summary of edits:
#  1. Split the "file-"BST.h"-and" run in the main.cpp bullet so a
#     gramStart/gramEnd proofErr pair wraps the quoted file name token.
#  2. Fill in the previously-empty bullet under "Changes to BSTNode.h" with
#     a new sentence, and append six new bulleted paragraphs describing the
#     isLeaf() rework (new ListParagraph items at ilvl 1/2/3, numId 4).
#  3. Split the "text book" integrity-statement sentence so a
#     gramStart/gramEnd proofErr pair wraps "text book".
#  4. Since list level 3 (ilvl=3) of numId 4 is now actually used in the
#     document body, numbering.xml's lvl 3 definition must stop being
#     "tentative" (w:tentative="1" removed), matching what Word itself
#     does once a previously-unused list level becomes used.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: "In main.cpp, I will inherit the appropriate file-"BST.h"-and..."
# ---------------------------------------------------------------------
$mainCppPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "In main.cpp, I will inherit*") {
        $mainCppPara = $cand
        break
    }
}

$mainCppXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>In main.cpp, I will inherit the appropriate file</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>-&#8220;</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/><w:r><w:t>BST.h</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&#8221;-and write my main() function that will handle creation and manipulation of the BST object through the BST implementation.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$mainCppPara.Range.InsertXML($mainCppXml)

Write-Output "edit 1 done"

# ---------------------------------------------------------------------
# Edit 2: fill the empty bullet after "...no left or right child node
# exists." with new text, then append six new bulleted paragraphs about
# the isLeaf() rework.
# ---------------------------------------------------------------------
$emptyPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    $prevText = ""
    if ($i -gt 1) { $prevText = $d.Paragraphs.Item($i - 1).Range.Text }
    if (($cand.Range.Text -eq "`r") -and ($prevText -like "*no left or right child node exists.*")) {
        $emptyPara = $cand
        break
    }
}

$isLeafXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">So, create context Booleans indicating the type of pointer of the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>lc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>rc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> pointers.</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Create getter and setter methods for these Booleans.</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Integrate these new variables into the constructor w/ parameters method.</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Amend </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>isLeaf</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>) to incorporate a check for the context variable when deciding if the node has no children (because now, all of them will have at least one child).</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>If it has no children, it&#8217;s a leaf</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">If </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>all of</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> its children are threaded, it&#8217;s a leaf</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="4"/></w:numPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>Have to</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> count number of children and get the &#8220;true&#8221; thread context variables for them</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>If they&#8217;re equal, it&#8217;s a leaf</w:t></w:r></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$emptyPara.Range.InsertXML($isLeafXml)

Write-Output "edit 2 done"

# ---------------------------------------------------------------------
# Edit 3: "If any source code or documentation ... such as a text book or
# course notes ..." -- wrap "text book" with a gramStart/gramEnd proofErr
# pair (split the run into three runs).
# ---------------------------------------------------------------------
$citationPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "If any source code or documentation*") {
        $citationPara = $cand
        break
    }
}

$citationXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">If any source code or documentation used in my program was obtained from another source, such as a </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:szCs w:val="24"/></w:rPr><w:t>text book</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> or course notes, that has been clearly noted with a proper citation in the comments of my program.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$citationPara.Range.InsertXML($citationXml)

Write-Output "edit 3 done"

# ---------------------------------------------------------------------
# Edit 4: numbering.xml -- list level 3 (ilvl=3) of the abstract numbering
# definition backing numId 4 is now genuinely used in the document body
# (see the new "Have to count..."/"If they're equal..." bullets above),
# so it is no longer "tentative" and the w:tentative="1" marker on that
# <w:lvl w:ilvl="3" w:tplc="0409000F"> must be dropped, exactly as real
# Word does once a previously-unused list level gets used.
# ---------------------------------------------------------------------
$fullXml = $d.Content.WordOpenXML
$needle = '<w:lvl w:ilvl="3" w:tplc="0409000F" w:tentative="1">'
$replacement = '<w:lvl w:ilvl="3" w:tplc="0409000F">'
$idx = $fullXml.IndexOf($needle)
if ($idx -ge 0) {
    $fullXml = $fullXml.Substring(0, $idx) + $replacement + $fullXml.Substring($idx + $needle.Length)
    $d.Content.InsertXML($fullXml)
    Write-Output "edit 4 done"
} else {
    Write-Output "edit 4 SKIPPED - needle not found"
}
